# Add data for 2022-10-28 (one more day of carjacking counts, through
# October 20 instead of October 19).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the sheet and update the "through October NN" label --------
$ws.Name = "Through 2022-10-20"
$ws.Range("B1").Value = "October 2022 (through October 20)"

# --- Bump existing counts by the new day's incidents --------------------
$ws.Range("V2").Value  = 13   # Garfield Park      / October 2020
$ws.Range("AP3").Value = 3    # Humboldt Park       / October 2018
$ws.Range("V4").Value  = 3    # South Shore         / October 2020
$ws.Range("L6").Value  = 11   # Austin              / October 2021
$ws.Range("AF6").Value = 4    # Austin              / October 2019
$ws.Range("AZ6").Value = 6    # Austin              / October 2017
$ws.Range("AZ7").Value = 4    # Englewood           / October 2017
$ws.Range("BJ7").Value = 4    # Englewood           / October 2016
$ws.Range("AF8").Value = 2    # Washington Heights  / October 2019
$ws.Range("B50").Value = 3    # Edgewater           / October 2022
$ws.Range("L58").Value = 2    # Archer Heights      / October 2021

# --- Fill in cells that were previously blank ---------------------------
$ws.Range("AF14").Value = 1   # Little Village       / October 2019
$ws.Range("BT16").Value = 1   # Little Italy, UIC    / October 2015
$ws.Range("V24").Value  = 1   # Auburn Gresham       / October 2020
$ws.Range("BJ25").Value = 1   # Hyde Park            / October 2016
$ws.Range("L60").Value  = 1   # Avalon Park          / October 2021
$ws.Range("B78").Value  = 1   # Jefferson Park       / October 2022
$ws.Range("B95").Value  = 1   # United Center        / October 2022
